$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append newly downloaded/analyzed participants sub_026, sub_027, sub_028
# following the existing pattern: column A = participant id, column B = boolean (FALSE)
$newParticipants = @("sub_026", "sub_027", "sub_028")

$row = 27
foreach ($p in $newParticipants) {
    $ws.Cells.Item($row, 1).Value = $p
    $ws.Cells.Item($row, 2).Value = $false
    $row++
}

# Update the selected cell to reflect the new last entry (B29)
$ws.Range("B29").Select()
